$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held two duplicate blocks of 7 shifts each:
#   rows 2-8  -> Sergio Soto   (RUN 17459567-4, terminal "LO MARCOLETA")
#   rows 9-15 -> Roberto Cordova (RUN 16808962-7, terminal "LO ESPEJO")
# The correction removes the Sergio Soto block entirely and keeps only
# Roberto Cordova's shifts (now starting at row 2), with the terminal name
# corrected to "Lo Marcoleta" and the week shifted forward by 7 days.

# Drop the Sergio Soto block; Roberto Cordova's rows shift up to 2:8.
$ws.Rows("2:8").Delete()

# Re-affirm / correct Roberto Cordova's data on every remaining row.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "16808962-7"
    $ws.Cells.Item($r, 2).Value2 = "Roberto "
    $ws.Cells.Item($r, 3).Value2 = "Cordova"
    $ws.Cells.Item($r, 4).Value2 = "Lo Marcoleta"
}

# Shift the week of dates forward by 7 days and update the recorded times.
$ws.Range("G2").Value2 = 44543
$ws.Range("G3").Value2 = 44544
$ws.Range("G4").Value2 = 44545
$ws.Range("G5").Value2 = 44546
$ws.Range("G6").Value2 = 44547
$ws.Range("G7").Value2 = 44548
$ws.Range("G8").Value2 = 44549

$ws.Range("H2").Value2 = 0.21180555555555555
$ws.Range("H3").Value2 = 0.21527777777777779
$ws.Range("H4").Value2 = 0.21875
$ws.Range("H5").Value2 = 0.22222222222222199
$ws.Range("H6").Value2 = 0.225694444444444
$ws.Range("H7").Value2 = 0.22916666666666699
$ws.Range("H8").Value2 = 0.23263888888888901

# Leave the selection where the author left it after the cleanup.
$ws.Range("I1:XFD1048576").Select()
